# The document contains five inline pictures. The first run (the title-page
# image) already carries <w:noProof/> on its run properties; the other four
# runs (the screenshots inserted later under "Praktikum" / "Tugas") are
# missing it. Word stamps <w:noProof/> on a run's rPr whenever the picture
# it holds was inserted/updated while the spell/grammar checker shouldn't
# re-evaluate that run - restore that flag on the four runs that still lack
# it.
$d = $word.ActiveDocument
$shapes = $d.InlineShapes

for ($i = 1; $i -le $shapes.Count; $i++) {
    $shape = $shapes.Item($i)
    $rng = $shape.Range
    if (-not $rng.NoProofing) {
        $rng.NoProofing = $true
    }
}
